$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.035.78"
$ws.Range("E2").Value = "  +0.47%  "
$ws.Range("D3").Value = "1.636.77"
$ws.Range("E3").Value = "  -0.02%  "
$ws.Range("E4").Value = "  +0.50%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.28%  "
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("E7").Value = "  +0.51%  "
$ws.Range("E8").Value = "  -1.80%  "
$ws.Range("E9").Value = "  -1.45%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.70"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.31%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0793"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.03%  "
$ws.Range("D12").Value = "1.697.68"
$ws.Range("E12").Value = "  +4.50%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "1.865.82"
$ws.Range("E13").Value = "  +0.10%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.21"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.53%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.533"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.97%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "0.0₃0749"
$ws.Range("E16").Value = "  -1.87%  "
$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.27"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.93%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "26.051.41"
$ws.Range("E18").Value = "  +0.39%  "
$ws.Range("B19").Value = "Dai"
$ws.Range("C19").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.01"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.42%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "191.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.04%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.28"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.83%  "
$ws.Range("B22").Value = "Avalanche"
$ws.Range("C22").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.64"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.80%  "
$ws.Range("B23").Value = "Chainlink"
$ws.Range("C23").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.16"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.71%  "
$ws.Range("B24").Value = "Stellar"
$ws.Range("C24").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.131"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.36%  "
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.63"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.27%  "
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.79"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("B27").Value = "BinanceUSD"
$ws.Range("C27").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.01"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.68%  "
$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.78"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.56%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.29"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.64%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.24"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.30%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0485"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.08%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.17"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.36%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.18"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.12%  "
$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.51"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.56%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.44"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.71%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.880"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.29%  "
$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D37").Value = "1.132.49"
$ws.Range("E37").Value = "  -0.10%  "
$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.45"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.528"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.93%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0156"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.77%  "
$ws.Range("B41").Value = "Quant"
$ws.Range("C41").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "98.93"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.41%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.788"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.04%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.32"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.97%  "
$ws.Range("B44").Value = "BabyDogeCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D44").Value = "0.0₆0113"
$ws.Range("E44").Value = "  -1.10%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "55.58"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.74%  "
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0525"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.65%  "
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.49"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.40%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.414"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.03%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.57"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.16%  "
$ws.Range("B50").Value = "USDD"
$ws.Range("C50").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.31%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0928"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.08%  "
